$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# for every data row (2 through 388) as part of an automatic file update.
$newValue = 46083

for ($row = 2; $row -le 388; $row++) {
    $ws.Cells.Item($row, 3).Value = $newValue
}
